$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values per row (column letter -> new value)
$updates = @{
    2 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "M" = 73.870127; "N" = 147.740254; "O" = 0.5038991021785622; "P" = 0.4272874344759938; "Q" = 1252.739365196534; "R" = 5010.957460786138; "S" = 0.3770848081571034; "T" = 0.2852693805155069 }
    3 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "O" = 0.0287796403411505; "P" = 0.0366060744105703; "Q" = 71.54882438901733; "R" = 429.292946334104; "S" = 0.02153678208584606; "T" = 0.02443926810769493 }
    4 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "M" = 14.18890566666667; "N" = 42.566717; "O" = 0.09678847345598375; "P" = 0.1231094627805071; "Q" = 240.6250185232999; "R" = 1443.750111139799; "S" = 0.07243010115948817; "T" = 0.08219141811661498 }
    5 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "M" = 20.1578925; "N" = 40.31578500000001; "O" = 0.1375054348094186; "P" = 0.1165994228054851; "Q" = 341.8511173555988; "R" = 1367.404469422395; "S" = 0.1028999858929986; "T" = 0.07784512819333834 }
    6 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "M" = 25.26300833333334; "N" = 75.78902500000001; "O" = 0.1723295699423422; "P" = 0.2191934640486469; "Q" = 428.4271099527792; "R" = 2570.562659716675; "S" = 0.1289600686735831; "T" = 0.1463398608454015 }
    7 = @{ "G" = 16.9586735; "H" = 33.917347; "I" = 0.7483339552041497; "J" = 0.6676287611063231; "M" = 8.898116; "N" = 26.694348; "O" = 0.06069777927254273; "P" = 0.0772041414787968; "Q" = 150.900244009126; "R" = 905.4014640547559; "S" = 0.04542220923513036; "T" = 0.05154370532776639 }
    8 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "M" = 73.870127; "N" = 147.740254; "O" = 0.5038991021785622; "P" = 0.4272874344759938; "Q" = 147.159560931653; "R" = 882.957365589918; "S" = 0.04429623299471411; "T" = 0.05026598263397029 }
    9 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "O" = 0.0287796403411505; "P" = 0.0366060744105703; "Q" = 8.404855690482668; "R" = 75.643701214344; "S" = 0.002529930393890494; "T" = 0.00430632907067837 }
    10 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "M" = 14.18890566666667; "N" = 42.566717; "O" = 0.09678847345598375; "P" = 0.1231094627805071; "Q" = 28.26627234588767; "R" = 254.396451112989; "S" = 0.008508379461032823; "T" = 0.0144825651748719 }
    11 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "M" = 20.1578925; "N" = 40.31578500000001; "O" = 0.1375054348094186; "P" = 0.1165994228054851; "Q" = 40.15732380705751; "R" = 240.9439428423451; "S" = 0.01208768333188869; "T" = 0.01371672576578135 }
    12 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "M" = 25.26300833333334; "N" = 75.78902500000001; "O" = 0.1723295699423422; "P" = 0.2191934640486469; "Q" = 50.32742415815834; "R" = 452.9468174234251; "S" = 0.01514896682498918; "T" = 0.02578586208803689 }
    13 = @{ "G" = 1.992139; "H" = 5.976417000000001; "I" = 0.0879069496317881; "J" = 0.117639739852435; "M" = 8.898116; "N" = 26.694348; "O" = 0.06069777927254273; "P" = 0.0772041414787968; "Q" = 17.726283910124; "R" = 159.536555191116; "S" = 0.005335756625272806; "T" = 0.009082275119096246 }
    14 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "M" = 73.870127; "N" = 147.740254; "O" = 0.5038991021785622; "P" = 0.4272874344759938; "Q" = 16.5719134859895; "R" = 66.28765394395799; "S" = 0.004988281673961843; "T" = 0.003773697566662749 }
    15 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "O" = 0.0287796403411505; "P" = 0.0366060744105703; "Q" = 0.9464865244439999; "R" = 5.678919146664001; "S" = 0.0002849001950515514; "T" = 0.0003232958490755847 }
    16 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "M" = 14.18890566666667; "N" = 42.566717; "O" = 0.09678847345598375; "P" = 0.1231094627805071; "Q" = 3.1831178139015; "R" = 19.098706883409; "S" = 0.0009581445299343631; "T" = 0.001087272507083424 }
    17 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "M" = 20.1578925; "N" = 40.31578500000001; "O" = 0.1375054348094186; "P" = 0.1165994228054851; "Q" = 4.52219136661125; "R" = 18.088765466445; "S" = 0.00136121663556153; "T" = 0.001029777434608976 }
    18 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "M" = 25.26300833333334; "N" = 75.78902500000001; "O" = 0.1723295699423422; "P" = 0.2191934640486469; "Q" = 5.6674653949875; "R" = 34.00479236992501; "S" = 0.001705953497254878; "T" = 0.001935862782679677 }
    19 = @{ "G" = 0.2243385; "H" = 0.448677; "I" = 0.009899366068316968; "J" = 0.008831754135926424; "M" = 8.898116; "N" = 26.694348; "O" = 0.06069777927254273; "P" = 0.0772041414787968; "Q" = 1.996189996266; "R" = 11.977139977596; "S" = 0.0006008695365528025; "T" = 0.0006818479958160124 }
    20 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "M" = 73.870127; "N" = 147.740254; "O" = 0.5038991021785622; "P" = 0.4272874344759938; "Q" = 216.6205031679013; "R" = 1299.723019007408; "S" = 0.0652045454539926; "T" = 0.07399208302515099 }
    21 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "O" = 0.0287796403411505; "P" = 0.0366060744105703; "Q" = 12.37204064214044; "R" = 111.348365779264; "S" = 0.0037240855533597; "T" = 0.006338964075396763 }
    22 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "M" = 14.18890566666667; "N" = 42.566717; "O" = 0.09678847345598375; "P" = 0.1231094627805071; "Q" = 41.60826588148711; "R" = 374.474392933384; "S" = 0.01252442877869399; "T" = 0.0213184963007584 }
    23 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "M" = 20.1578925; "N" = 40.31578500000001; "O" = 0.1375054348094186; "P" = 0.1165994228054851; "Q" = 59.11202530022; "R" = 354.6721518013201; "S" = 0.01779320370970727; "T" = 0.02019117221054823 }
    24 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "M" = 25.26300833333334; "N" = 75.78902500000001; "O" = 0.1723295699423422; "P" = 0.2191934640486469; "Q" = 74.08252562908889; "R" = 666.7427306618001; "S" = 0.02229944690870002; "T" = 0.03795707451670717 }
    25 = @{ "G" = 2.932450666666667; "H" = 8.797352; "I" = 0.1294000032389156; "J" = 0.1731669996705884; "M" = 8.898116; "N" = 26.694348; "O" = 0.06069777927254273; "P" = 0.0772041414787968; "Q" = 26.09328619627733; "R" = 234.839575766496; "S" = 0.007854292834462014; "T" = 0.01336920954202687 }
    26 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "M" = 73.870127; "N" = 147.740254; "O" = 0.5038991021785622; "P" = 0.4272874344759938; "Q" = 40.946506876608; "R" = 245.679041259648; "S" = 0.0123252338987902; "T" = 0.01398629073470288 }
    27 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "O" = 0.0287796403411505; "P" = 0.0366060744105703; "Q" = 2.338614488576; "R" = 21.047530397184; "S" = 0.0007039421130026952; "T" = 0.001198217307724663 }
    28 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "M" = 14.18890566666667; "N" = 42.566717; "O" = 0.09678847345598375; "P" = 0.1231094627805071; "Q" = 7.864967166656001; "R" = 70.78470449990401; "S" = 0.002367419526834392; "T" = 0.004029710681178468 }
    29 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "M" = 20.1578925; "N" = 40.31578500000001; "O" = 0.1375054348094186; "P" = 0.1165994228054851; "Q" = 11.17360044432; "R" = 67.04160266592001; "S" = 0.003363345239262534; "T" = 0.00381661920120818 }
    30 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "M" = 25.26300833333334; "N" = 75.78902500000001; "O" = 0.1723295699423422; "P" = 0.2191934640486469; "Q" = 14.0033865712; "R" = 126.0304791408; "S" = 0.004215134037815034; "T" = 0.007174803815821687 }
    31 = @{ "G" = 0.554304; "H" = 1.662912; "I" = 0.02445972585682961; "J" = 0.03273274523472717; "M" = 8.898116; "N" = 26.694348; "O" = 0.06069777927254273; "P" = 0.0772041414787968; "Q" = 4.932261291264; "R" = 44.390351621376; "S" = 0.001484651041124749; "T" = 0.002527103494091288 }
}

foreach ($rowNum in $updates.Keys) {
    $rowUpdates = $updates[$rowNum]
    foreach ($col in $rowUpdates.Keys) {
        $ws.Range("$col$rowNum").Value = $rowUpdates[$col]
    }
}
